# RegisterCommand UML diagram - resolved discrepancies
#  1. Bump the cached "datetimeFigureOut" date field text (master + all
#     layouts) from 15-10-2018 to 10-11-2018.
#  2. Replace the single-line "Command" superclass box with a taller
#     box that shows the UML "{abstract}" stereotype above the
#     "Command" class name.

$p = $ppt.ActivePresentation

# --- 1. Update the cached date/time field text everywhere it appears ---

function Update-DateShapes($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.HasText) {
                $t = $sh.TextFrame.TextRange.Text
                if ($t -eq "15-10-2018") {
                    $sh.TextFrame.TextRange.Text = "10-11-2018"
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShapes $master.Shapes

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DateShapes $layouts.Item($i).Shapes
}

# --- 2. Rebuild the "Command" abstract-class box on slide 1 ---

$s = $p.Slides.Item(1)

# Locate the existing "Command" rectangle (id 10 in the source XML) by
# its text so the script is resilient to ordering.
$cmdShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq "Command") {
            $cmdShape = $sh
        }
    }
}

# Duplicate it so the new shape inherits the exact same style (line /
# fill / effect / font refs) as the original.
$newShape = $cmdShape.Duplicate().Item(1)

# Reposition / resize: grow the box upward so it can hold two lines.
$newShape.Left = 330.07543307086615
$newShape.Top = 85.1615748031496
$newShape.Width = 125.12622047244095
$newShape.Height = 46.477798

# The new box has no explicit outline override (inherits from style).
$newShape.Line.Visible = $true

# Insert the "{abstract}" stereotype as its own paragraph above the
# existing "Command" run, preserving that run's formatting/endParaRPr.
$newShape.TextFrame.TextRange.InsertBefore("{abstract}" + [char]13)
$firstPara = $newShape.TextFrame.TextRange.Paragraphs(1, 1)
$firstPara.Font.Size = 14
$firstPara.Font.Bold = $true
$firstPara.Font.Italic = $true
$firstPara.ParagraphFormat.Alignment = 2

# Move the new shape to sit right after the "1" multiplicity label
# (and before the "User" class box), matching its place in the XML.
for ($k = 0; $k -lt 4; $k++) {
    $newShape.ZOrder(3)
}

# Remove the old single-line "Command" box.
$cmdShape.Delete()
